$wb = $excel.ActiveWorkbook

# Update the SiPM_Strip_ID value on the Instructions sheet.
# This value is referenced by formulas (=Instructions!$B$3) throughout
# the other worksheets, so updating it here propagates everywhere.
$wsInstructions = $wb.Worksheets.Item("Instructions")
$wsInstructions.Range("B3").Value = "Test-010"

# Update the active cell / selection on the SiPM-item-manifest sheet (A1 -> J2).
# Activate it first so Range.Select() applies to this sheet's view.
$wsManifest = $wb.Worksheets.Item("SiPM-item-manifest")
$wsManifest.Activate() | Out-Null
$wsManifest.Range("J2").Select() | Out-Null

# Update the active cell / selection on the Instructions sheet (B3 -> B4).
# Re-activate Instructions last so it remains the selected tab, matching
# the original workbook's tabSelected state.
$wsInstructions.Activate() | Out-Null
$wsInstructions.Range("B4").Select() | Out-Null
